# Add the new COVID totals rows (564-569) to the bottom of the existing
# table on Sheet1, extending the sheet dimension from A1:H563 to A1:H569.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "date" column holds text values such as "2022-02-28" which Excel
# would otherwise auto-detect and convert into a date serial number.
# Force the new date cells to be stored as plain text first so the
# values round-trip as strings, matching the rest of column A.
$ws.Range("A564:A569").NumberFormat = "@"

$rows = @(
    @("2022-02-28", "overview", "K02000001", "United Kingdom", 18886701, 82451, 138, 161361),
    @("2022-03-01", "overview", "K02000001", "United Kingdom", 18985568, 39000, 194, 161630),
    @("2022-03-02", "overview", "K02000001", "United Kingdom", 19029321, 44017, 74, 161704),
    @("2022-03-03", "overview", "K02000001", "United Kingdom", 19074696, 45656, 194, 161898),
    @("2022-03-04", "overview", "K02000001", "United Kingdom", 19119181, 44740, 110, 162008),
    @("2022-03-07", "overview", "K02000001", "United Kingdom", 19245301, 126604, 139, 162147)
)

$startRow = 564
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
